# Revert "updating testcase of adding 20 items into cart"
# Re-insert the "Go to My Account" CLICK step above the QuickOrder step
# (old row 9), shifting every following row down by one, and restore the
# previously-active sheet/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new blank row at row 9 (pushes existing rows 9..39 down to 10..40)
$ws1.Rows.Item(9).Insert()

# Populate the new row 9 with the same step as row 39 (CLICK / MyaccountSection / xpath)
$ws1.Range("A9").Value = ""
$ws1.Range("B9").Value = "CLICK"
$ws1.Range("C9").Value = "MyaccountSection"
$ws1.Range("D9").Value = "xpath"
$ws1.Range("E9").Value = ""

# Match formatting of the equivalent row further down (row 39, the other
# MyaccountSection/CLICK row) so the re-inserted row looks the same as its
# sibling rows.
$ws1.Range("A9:E9").Borders.LineStyle = 1
$ws1.Range("B9,C9,D9,E9").Font.Name = "Calibri"

# Restore the view state: selection + scroll position on sheet 1
$ws1.Range("E37:E38").Select()
$ws1.Application.ActiveWindow.ScrollRow = 22

# Sheet 2 (Testdata) becomes the active sheet with its own selection restored
$ws2.Activate()
$ws2.Range("B27:B28").Select()
$ws2.Application.ActiveWindow.ScrollRow = 19
